$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '86.979.27'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +9.45%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.314.45'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +5.12%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.15%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.59'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.05%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '639.32'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.16%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.323'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +21.44%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.997'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.20%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.609'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.15%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.322.94'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.44%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.603'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.44%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000276'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +10.32%  '

# Row 13
$ws.Range("E13").Value = '  +1.44%  '

# Row 14
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.37'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +9.46%  '

# Row 15
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.912.21'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.56%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.40'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.61%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '86.636.30'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +8.85%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.311.29'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.64%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.65'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.48%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.19'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +8.21%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '448.60'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.58%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.16'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.71%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.26'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.77%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.42'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +9.68%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.39'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +15.26%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.28'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +13.35%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.440.55'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.22%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '78.54'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.34%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000131'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +8.15%  '

# Row 30
$ws.Range("E30").Value = '  -0.03%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.173'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +44.73%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '609.47'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +11.05%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.26'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.43%  '

# Row 34
$ws.Range("E34").Value = '  +0.10%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.56'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.02%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.05'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.98%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.152'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.15%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '23.44'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.39%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.55'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +16.53%  '

# Row 40
$ws.Range("B40").Value = 'PolygonEcosystemToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.418'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.10%  '

# Row 41
$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.997'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.11%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '21.34'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.73%  '

# Row 43
$ws.Range("E43").Value = '  +16.66%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.07'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +15.01%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '158.15'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.30%  '

# Row 46
$ws.Range("E46").Value = '  +0.05%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '188.68'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.43%  '

# Row 48
$ws.Range("B48").Value = 'OKB'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '45.64'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.08%  '

# Row 49
$ws.Range("B49").Value = 'ImmutableX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.37'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.23%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.785'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.26%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '26.43'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.74%  '

Write-Host "Applied all cryptos list updates"